$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell E1, copying the style from D1 (bold header style)
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "QUANTIDADE"

# Fill in the QUANTIDADE values for rows 2-47
$quantidades = @(
    753, 220, 1040, 4173, 12547, 29668, 58364, 96056, 139131, 178581, 207272, 220086, 219114, 204504, 183746, 159655, 135809, 113998, 95703, 79025, 65788, 54547, 44896, 36990, 30709, 25307, 20900, 17479, 14373, 11867, 9995, 8149, 6580, 5448, 4350, 3479, 2515, 1794, 1188, 735, 392, 187, 58, 22, 0, 0
)

for ($i = 0; $i -lt $quantidades.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $quantidades[$i]
}

$ws.Range("A1").Select()
